{"js": "// Replace the date line and each \"A\u00f7B=\" equation with its updated value.\n// All old strings are unique within the document, so a plain search/replace\n// (with wildcards off, match case on) is unambiguous for every entry.\nconst replacements = [\n  [\"2024-09-23 Monday\", \"2024-09-24 Tuesday\"],\n  [\"423\u00f78=\", \"638\u00f73=\"],\n  [\"430\u00f75=\", \"227\u00f75=\"],\n  [\"526\u00f76=\", \"659\u00f75=\"],\n  [\"809\u00f76=\", \"387\u00f76=\"],\n  [\"162\u00f76=\", \"957\u00f76=\"],\n  [\"269\u00f76=\", \"256\u00f72=\"],\n  [\"728\u00f74=\", \"746\u00f79=\"],\n  [\"626\u00f74=\", \"801\u00f72=\"],\n  [\"944\u00f74=\", \"120\u00f72=\"],\n  [\"817\u00f72=\", \"586\u00f77=\"],\n  [\"925\u00f76=\", \"898\u00f73=\"],\n  [\"196\u00f74=\", \"730\u00f77=\"],\n  [\"381\u00f78=\", \"801\u00f77=\"],\n  [\"528\u00f73=\", \"670\u00f76=\"],\n  [\"336\u00f78=\", \"724\u00f77=\"],\n  [\"429\u00f76=\", \"948\u00f79=\"],\n  [\"941\u00f77=\", \"502\u00f72=\"],\n  [\"383\u00f77=\", \"184\u00f76=\"],\n  [\"384\u00f74=\", \"676\u00f76=\"],\n  [\"199\u00f75=\", \"390\u00f73=\"],\n  [\"763\u00f74=\", \"733\u00f77=\"],\n  [\"218\u00f74=\", \"771\u00f75=\"],\n  [\"165\u00f76=\", \"217\u00f74=\"],\n  [\"367\u00f72=\", \"733\u00f76=\"],\n  [\"316\u00f77=\", \"108\u00f76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < found.items.length; i++) {\n    found.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each \"A\u00f7B=\" equation with its updated value.\n# Every \"old\" string below is unique in the document, so a plain\n# Find/Replace (wildcards off, match case on) unambiguously targets the\n# correct run each time.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old = \"2024-09-23 Monday\"; New = \"2024-09-24 Tuesday\"},\n    @{Old = \"423\u00f78=\"; New = \"638\u00f73=\"},\n    @{Old = \"430\u00f75=\"; New = \"227\u00f75=\"},\n    @{Old = \"526\u00f76=\"; New = \"659\u00f75=\"},\n    @{Old = \"809\u00f76=\"; New = \"387\u00f76=\"},\n    @{Old = \"162\u00f76=\"; New = \"957\u00f76=\"},\n    @{Old = \"269\u00f76=\"; New = \"256\u00f72=\"},\n    @{Old = \"728\u00f74=\"; New = \"746\u00f79=\"},\n    @{Old = \"626\u00f74=\"; New = \"801\u00f72=\"},\n    @{Old = \"944\u00f74=\"; New = \"120\u00f72=\"},\n    @{Old = \"817\u00f72=\"; New = \"586\u00f77=\"},\n    @{Old = \"925\u00f76=\"; New = \"898\u00f73=\"},\n    @{Old = \"196\u00f74=\"; New = \"730\u00f77=\"},\n    @{Old = \"381\u00f78=\"; New = \"801\u00f77=\"},\n    @{Old = \"528\u00f73=\"; New = \"670\u00f76=\"},\n    @{Old = \"336\u00f78=\"; New = \"724\u00f77=\"},\n    @{Old = \"429\u00f76=\"; New = \"948\u00f79=\"},\n    @{Old = \"941\u00f77=\"; New = \"502\u00f72=\"},\n    @{Old = \"383\u00f77=\"; New = \"184\u00f76=\"},\n    @{Old = \"384\u00f74=\"; New = \"676\u00f76=\"},\n    @{Old = \"199\u00f75=\"; New = \"390\u00f73=\"},\n    @{Old = \"763\u00f74=\"; New = \"733\u00f77=\"},\n    @{Old = \"218\u00f74=\"; New = \"771\u00f75=\"},\n    @{Old = \"165\u00f76=\"; New = \"217\u00f74=\"},\n    @{Old = \"367\u00f72=\"; New = \"733\u00f76=\"},\n    @{Old = \"316\u00f77=\"; New = \"108\u00f76=\"}\n)\n\nforeach ($r in $replacements) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute(\n        $r.Old, $true, $false, $false, $false, $false, $true, 1, $false,\n        $r.New, 2\n    ) | Out-Null\n}\n"}
